$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Rows 17 and 18 swap coin identity (ShibaInu <-> WrappedEther ranking)
# and receive brand-new Price / Volume(1h) values
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D17' '2.621.17'
$ws.Range('E17').Value = '  +6.97%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D18' '0.0000131'
$ws.Range('E18').Value = '  -1.70%  '

# Remaining Price / Volume(1h) updates
Set-TextValue 'D2' '56.020.19'
$ws.Range('E2').Value = '  -3.28%  '
Set-TextValue 'D3' '2.357.25'
$ws.Range('E3').Value = '  -3.72%  '
$ws.Range('E4').Value = '  +0.15%  '
Set-TextValue 'D5' '499.26'
$ws.Range('E5').Value = '  -1.95%  '
Set-TextValue 'D6' '128.65'
$ws.Range('E6').Value = '  -3.17%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  -0.04%  '
Set-TextValue 'D8' '0.543'
$ws.Range('E8').Value = '  -2.52%  '
Set-TextValue 'D9' '2.355.22'
$ws.Range('E9').Value = '  -3.80%  '
Set-TextValue 'D10' '0.0973'
$ws.Range('E10').Value = '  -0.64%  '
$ws.Range('E11').Value = '  +0.23%  '
Set-TextValue 'D12' '4.72'
$ws.Range('E12').Value = '  +2.94%  '
$ws.Range('E13').Value = '  -0.32%  '
Set-TextValue 'D14' '2.773.55'
$ws.Range('E14').Value = '  -3.83%  '
Set-TextValue 'D15' '55.953.22'
$ws.Range('E15').Value = '  -2.78%  '
Set-TextValue 'D16' '21.29'
$ws.Range('E16').Value = '  -2.57%  '
Set-TextValue 'D19' '9.94'
$ws.Range('E19').Value = '  -3.31%  '
$ws.Range('E20').Value = '  -2.42%  '
Set-TextValue 'D21' '305.06'
$ws.Range('E21').Value = '  -2.82%  '
$ws.Range('E22').Value = '  -2.95%  '
$ws.Range('E23').Value = '  -0.07%  '
Set-TextValue 'D24' '65.02'
$ws.Range('E24').Value = '  -0.55%  '
$ws.Range('E25').Value = '  +0.85%  '
Set-TextValue 'D26' '0.367'
$ws.Range('E26').Value = '  -3.66%  '
$ws.Range('E27').Value = '  -6.46%  '
Set-TextValue 'D29' '171.06'
$ws.Range('E29').Value = '  -1.45%  '
Set-TextValue 'D30' '0.0₃0705'
$ws.Range('E30').Value = '  -3.83%  '
Set-TextValue 'D31' '1.62'
$ws.Range('E31').Value = '  -3.80%  '
$ws.Range('E32').Value = '  +0.10%  '
Set-TextValue 'D33' '0.998'
$ws.Range('E33').Value = '  -0.08%  '
Set-TextValue 'D34' '1.08'
$ws.Range('E34').Value = '  -5.13%  '
Set-TextValue 'D35' '5.71'
$ws.Range('E35').Value = '  -7.38%  '
Set-TextValue 'D36' '17.52'
$ws.Range('E36').Value = '  -2.69%  '
$ws.Range('E37').Value = '  -5.98%  '
$ws.Range('E38').Value = '  -2.96%  '
Set-TextValue 'D39' '35.95'
$ws.Range('E39').Value = '  -1.97%  '
Set-TextValue 'D40' '0.783'
$ws.Range('E40').Value = '  -3.57%  '
$ws.Range('E41').Value = '  -6.13%  '
Set-TextValue 'D42' '128.89'
$ws.Range('E42').Value = '  -5.16%  '
$ws.Range('E43').Value = '  -1.88%  '
Set-TextValue 'D44' '4.67'
$ws.Range('E44').Value = '  -5.27%  '
$ws.Range('E45').Value = '  -2.29%  '
Set-TextValue 'D46' '0.0900'
$ws.Range('E46').Value = '  -1.90%  '
Set-TextValue 'D47' '237.81'
$ws.Range('E47').Value = '  -7.12%  '
$ws.Range('E48').Value = '  -2.68%  '
$ws.Range('E49').Value = '  -3.61%  '
$ws.Range('E50').Value = '  -0.82%  '
Set-TextValue 'D51' '0.949'
$ws.Range('E51').Value = '  -0.74%  '
